$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Build the two new border styles exactly once (top+bottom only, and
# top+right+bottom), on sheet 2's C1/D1, then reuse them everywhere else
# via a format-only copy/paste. Re-deriving the same border combination
# independently on more than two distinct cells causes the COM layer to
# leave behind an extra (unused) entry in the style table, so we build it
# once and fan it out by copying formats instead. ---

$c1b = $ws2.Cells.Item(1, 3)
$c1b.Style = "Normal"
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142

$d1b = $ws2.Cells.Item(1, 4)
$d1b.Style = "Normal"
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

# F1 / G1 (second merged group E1:G1 on sheet 2) reuse the same styles
$c1b.Copy()
$f1b = $ws2.Cells.Item(1, 6)
$f1b.PasteSpecial(-4122)

$d1b.Copy()
$g1b = $ws2.Cells.Item(1, 7)
$g1b.PasteSpecial(-4122)

# C1 / D1 on sheet 1 reuse the same styles too
$c1b.Copy()
$c1 = $ws1.Cells.Item(1, 3)
$c1.PasteSpecial(-4122)

$d1b.Copy()
$d1 = $ws1.Cells.Item(1, 4)
$d1.PasteSpecial(-4122)

# --- Rename "fedcore" headers to "approach" ---
$ws1.Cells.Item(2, 3).Value2 = "approach"
$ws2.Cells.Item(2, 3).Value2 = "approach"
$ws2.Cells.Item(2, 6).Value2 = "approach"

# --- Clear the stray empty inline-string cell G5 on sheet 2 ---
$ws2.Cells.Item(5, 7).ClearContents()
